# Fruta / hortaliza, semanal
# Insert a new weekly record at row 89 (shifting the existing rows 89:201 down to 90:202)
# and populate it with the new week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 89; this shifts rows 89:201 -> 90:202
# and Excel naturally carries formatting (e.g. the date number format on column D) down
# from the row above into the newly inserted row.
$ws.Rows.Item(89).Insert()

# Populate the newly inserted row 89 with the new record's values.
$ws.Range("A89").Value = 10
$ws.Range("B89").Value = "Vega Modelo de Temuco"
$ws.Range("C89").Value = "La Araucanía"
$ws.Range("D89").Value = 44482
$ws.Range("E89").Value = 9
$ws.Range("F89").Value = 100112009
$ws.Range("G89").Value = "Acelga"
$ws.Range("H89").Value = "Sin especificar"
$ws.Range("I89").Value = "Primera"
$ws.Range("J89").Value = 30
$ws.Range("K89").Value = 8000
$ws.Range("L89").Value = 8000
$ws.Range("M89").Value = 8000
$ws.Range("N89").Value = "$/docena de atados (12 kilos)"
$ws.Range("O89").Value = "Provincia de Cautín"
$ws.Range("P89").Value = 667
$ws.Range("Q89").Value = 12
$ws.Range("R89").Value = "Hortaliza"
